$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# The "Förändrad" (changed) date column C for rows 2-6 was updated
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06).
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45175
}
